$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.262.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.096.99"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5262"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4382"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.06"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09375"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.175"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.75"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.558"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.128.84"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.870"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.41"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.11"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06730"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.409"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.283.21"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.323"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.015"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.79"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.518"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.137"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.673"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.261"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.877"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.10"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02630"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06768"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.66"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.348"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6969"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2214"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6795"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.32"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.308"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.628"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000351"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.59%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07296"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.44%  "
